$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(318294931, Shalev  Afanasenko: -7,6)"
$ws.Range("B1").Value = "(305487936, Avihai  Kipnis: -9,-4)"
$ws.Range("C1").Value = "(313227928, Aviv  Levi: 7,6)"
$ws.Range("D1").Value = "(205807308, Sariel  Basis: 0,4)"
$ws.Range("E1").Value = "(315891549, Raz  Halaby: -8,-2)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: 6,1)"
$ws.Range("G1").Value = "(313925141, Elad   Amer: 1,8)"

$ws.Range("A3").Value = "cost: 605.0358990268988"
$ws.Range("A4").Value = "time: 82.862271289557"
